# Update "想去人数" (F column) counts on the three sheets that list events:
# 展览 (Exhibition), 演出 (Performance) and 全部类型 (All types, which mirrors
# the other sheets). Values are refreshed counters bumped in the source data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览"
$ws1.Range("F3").Value = 3220
$ws1.Range("F4").Value = 1994
$ws1.Range("F6").Value = 100
$ws1.Range("F7").Value = 3075
$ws1.Range("F8").Value = 610
$ws1.Range("F10").Value = 38
$ws1.Range("F15").Value = 10115
$ws1.Range("F20").Value = 8001
$ws1.Range("F21").Value = 12609
$ws1.Range("F24").Value = 19
$ws1.Range("F25").Value = 268
$ws1.Range("F26").Value = 395
$ws1.Range("F27").Value = 592
$ws1.Range("F30").Value = 2820
$ws1.Range("F32").Value = 237
$ws1.Range("F33").Value = 7937
$ws1.Range("F34").Value = 1478
$ws1.Range("F35").Value = 219
$ws1.Range("F38").Value = 4619
$ws1.Range("F39").Value = 1390
$ws1.Range("F43").Value = 637

# Sheet "演出"
$ws2.Range("F6").Value = 1197

# Sheet "全部类型"
$ws4.Range("F5").Value = 3220
$ws4.Range("F7").Value = 1994
$ws4.Range("F11").Value = 3075
$ws4.Range("F13").Value = 610
$ws4.Range("F14").Value = 38
$ws4.Range("F19").Value = 10115
$ws4.Range("F23").Value = 8001
$ws4.Range("F24").Value = 12609
$ws4.Range("F27").Value = 19
$ws4.Range("F28").Value = 268
$ws4.Range("F30").Value = 592
$ws4.Range("F33").Value = 2820
$ws4.Range("F37").Value = 237
$ws4.Range("F38").Value = 7937
$ws4.Range("F39").Value = 219
$ws4.Range("F42").Value = 4619
$ws4.Range("F47").Value = 637
